$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values for rows 2-51
# D column values are stored as literal text (e.g. "1.826.10"), so force
# Text number format before assigning to avoid numeric reinterpretation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.115.37"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.10"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  -0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.62"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4568"
$ws.Range("E7").Value = "  +7.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3735"
$ws.Range("E8").Value = "  +1.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07324"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8602"
$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.04"
$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.696"
$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("E14").Value = "  +5.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.349"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07078"
$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008841"
$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("E19").Value = "  -0.41%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.129.80"
$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.193"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  +1.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.997"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.69"
$ws.Range("E25").Value = "  -0.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.219"
$ws.Range("E26").Value = "  +5.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.54"
$ws.Range("E27").Value = "  +1.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.265"
$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.50"
$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08857"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7729"
$ws.Range("E31").Value = "  +1.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.194"
$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("E33").Value = "  +6.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.464"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.104"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01969"
$ws.Range("E37").Value = "  +0.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05288"
$ws.Range("E38").Value = "  +0.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5364"
$ws.Range("E39").Value = "  +7.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.189"
$ws.Range("E40").Value = "  +2.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.883"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1717"
$ws.Range("E42").Value = "  +2.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5246"
$ws.Range("E43").Value = "  +11.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.620"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.72"
$ws.Range("E45").Value = "  +1.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.995"
$ws.Range("E46").Value = "  +11.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.14"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06492"
$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.677"
$ws.Range("E49").Value = "  +1.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9234"
$ws.Range("E51").Value = "  +1.35%  "
